# Scheduled data refresh: updates market-price columns (H currentAveragePrice,
# I currentAveragePriceNQ, J currentAveragePriceHQ, K LevePriceNQ, L LevePriceHQ,
# M LeveProfitNQ, N LeveProfitHQ) for specific leve rows across the Ixion_Profits
# sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR), reflecting newly-polled prices.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1263.6471
$ws.Range("J32").Value = 1581.8334
$ws.Range("L32").Value = 1581.8334
$ws.Range("N32").Value = -2233.8334
$ws.Range("H121").Value = 1645
$ws.Range("I121").Value = 892.5
$ws.Range("J121").Value = 2397.5
$ws.Range("K121").Value = 2677.5
$ws.Range("L121").Value = 7192.5
$ws.Range("M121").Value = -930.5
$ws.Range("N121").Value = -10686.5
$ws.Range("H124").Value = 23588.889
$ws.Range("J124").Value = 23588.889
$ws.Range("L124").Value = 23588.889
$ws.Range("N124").Value = -33408.889
$ws.Range("H136").Value = 35700
$ws.Range("J136").Value = 35700
$ws.Range("L136").Value = 35700
$ws.Range("N136").Value = -45900
$ws.Range("H139").Value = 33463.332
$ws.Range("J139").Value = 33463.332
$ws.Range("L139").Value = 33463.332
$ws.Range("N139").Value = -43743.332
$ws.Range("H140").Value = 64000
$ws.Range("J140").Value = 64000
$ws.Range("L140").Value = 64000
$ws.Range("N140").Value = -74360

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5860.1816
$ws.Range("I45").Value = 6434.6665
$ws.Range("J45").Value = 3275
$ws.Range("K45").Value = 6434.6665
$ws.Range("L45").Value = 3275
$ws.Range("M45").Value = -6057.6665
$ws.Range("N45").Value = -4029
$ws.Range("H95").Value = 43000
$ws.Range("J95").Value = 43000
$ws.Range("L95").Value = 43000
$ws.Range("N95").Value = -48492

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1250372.5
$ws.Range("I7").Value = 496.66666
$ws.Range("J7").Value = 5000000
$ws.Range("K7").Value = 496.66666
$ws.Range("L7").Value = 5000000
$ws.Range("M7").Value = -383.66666
$ws.Range("N7").Value = -5000226

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 228.45454
$ws.Range("I7").Value = 272.2
$ws.Range("J7").Value = 192
$ws.Range("K7").Value = 272.2
$ws.Range("L7").Value = 192
$ws.Range("M7").Value = -159.2
$ws.Range("N7").Value = -418
$ws.Range("H31").Value = 2758.717
$ws.Range("I31").Value = 1386.9231
$ws.Range("K31").Value = 1386.9231
$ws.Range("M31").Value = -1091.9231
$ws.Range("H34").Value = 2758.717
$ws.Range("I34").Value = 1386.9231
$ws.Range("K34").Value = 1386.9231
$ws.Range("M34").Value = -1184.9231

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 39285830
$ws.Range("I12").Value = 100000090
$ws.Range("J12").Value = 136.88235
$ws.Range("K12").Value = 300000270
$ws.Range("L12").Value = 410.64705
$ws.Range("M12").Value = -300000097
$ws.Range("N12").Value = -756.64705
$ws.Range("H132").Value = 2920.8
$ws.Range("I132").Value = 2002
$ws.Range("J132").Value = 3533.3333
$ws.Range("K132").Value = 18018
$ws.Range("L132").Value = 31799.9997
$ws.Range("M132").Value = -15488
$ws.Range("N132").Value = -36859.9997
$ws.Range("H134").Value = 7398.2856
$ws.Range("I134").Value = 5860.909
$ws.Range("K134").Value = 17582.727
$ws.Range("M134").Value = -12512.727
$ws.Range("H137").Value = 33370428
$ws.Range("I137").Value = 18643.334
$ws.Range("J137").Value = 83398104
$ws.Range("K137").Value = 55930.00199999999
$ws.Range("L137").Value = 250194312
$ws.Range("M137").Value = -50830.00199999999
$ws.Range("N137").Value = -250204512
$ws.Range("H139").Value = 4439.617
$ws.Range("I139").Value = 7797.3335
$ws.Range("J139").Value = 2865.6875
$ws.Range("K139").Value = 23392.0005
$ws.Range("L139").Value = 8597.0625
$ws.Range("M139").Value = -18252.0005
$ws.Range("N139").Value = -18877.0625
$ws.Range("H140").Value = 5954.6924
$ws.Range("I140").Value = 3117.5833
$ws.Range("J140").Value = 40000
$ws.Range("K140").Value = 9352.749899999999
$ws.Range("L140").Value = 120000
$ws.Range("M140").Value = -4172.749899999999
$ws.Range("N140").Value = -130360
$ws.Range("H141").Value = 11064.692
$ws.Range("I141").Value = 11593.444
$ws.Range("J141").Value = 9875
$ws.Range("K141").Value = 34780.33199999999
$ws.Range("L141").Value = 29625
$ws.Range("M141").Value = -29600.33199999999
$ws.Range("N141").Value = -39985

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H70").Value = 5742.8286
$ws.Range("I70").Value = 5667.6
$ws.Range("J70").Value = 5930.9
$ws.Range("K70").Value = 5667.6
$ws.Range("L70").Value = 5930.9
$ws.Range("M70").Value = -5397.6
$ws.Range("N70").Value = -6470.9
$ws.Range("H73").Value = 5742.8286
$ws.Range("I73").Value = 5667.6
$ws.Range("J73").Value = 5930.9
$ws.Range("K73").Value = 5667.6
$ws.Range("L73").Value = 5930.9
$ws.Range("M73").Value = -4731.6
$ws.Range("N73").Value = -7802.9
$ws.Range("H80").Value = 3048.3333
$ws.Range("I80").Value = 2575
$ws.Range("J80").Value = 3995
$ws.Range("K80").Value = 2575
$ws.Range("L80").Value = 3995
$ws.Range("M80").Value = -1577
$ws.Range("N80").Value = -5991
$ws.Range("H83").Value = 3048.3333
$ws.Range("I83").Value = 2575
$ws.Range("J83").Value = 3995
$ws.Range("K83").Value = 12875
$ws.Range("L83").Value = 19975
$ws.Range("M83").Value = -7883
$ws.Range("N83").Value = -29959

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1069.4286
$ws.Range("I100").Value = 917.2
$ws.Range("J100").Value = 1450
$ws.Range("K100").Value = 917.2
$ws.Range("L100").Value = 1450
$ws.Range("M100").Value = -376.2
$ws.Range("N100").Value = -2532
$ws.Range("H122").Value = 4792916
$ws.Range("I122").Value = 5497920.5
$ws.Range("J122").Value = 2501650
$ws.Range("K122").Value = 16493761.5
$ws.Range("L122").Value = 7504950
$ws.Range("M122").Value = -16491311.5
$ws.Range("N122").Value = -7509850

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 25057
$ws.Range("J64").Value = 25057
$ws.Range("L64").Value = 25057
$ws.Range("N64").Value = -25553
$ws.Range("H67").Value = 25057
$ws.Range("J67").Value = 25057
$ws.Range("L67").Value = 25057
$ws.Range("N67").Value = -26773
$ws.Range("H81").Value = 2102.1738
$ws.Range("I81").Value = 1930.0834
$ws.Range("J81").Value = 2289.9092
$ws.Range("K81").Value = 3860.1668
$ws.Range("L81").Value = 4579.8184
$ws.Range("M81").Value = -2799.1668
$ws.Range("N81").Value = -6701.8184
$ws.Range("H84").Value = 2102.1738
$ws.Range("I84").Value = 1930.0834
$ws.Range("J84").Value = 2289.9092
$ws.Range("K84").Value = 19300.834
$ws.Range("L84").Value = 22899.092
$ws.Range("M84").Value = -13996.834
$ws.Range("N84").Value = -33507.092
$ws.Range("H122").Value = 1261.2778
$ws.Range("I122").Value = 1249.875
$ws.Range("J122").Value = 1352.5
$ws.Range("K122").Value = 3749.625
$ws.Range("L122").Value = 4057.5
$ws.Range("M122").Value = -1299.625
$ws.Range("N122").Value = -8957.5
